$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 58. This shifts the existing rows 58..203 down to 59..204,
# which is exactly what the target diff shows (every row from 58 downward is
# pushed down by one, and a brand-new record appears at the top of that block).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted (currently blank) row 58 with the new record.
$ws.Cells.Item(58, 1).Value2 = 4
$ws.Cells.Item(58, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(58, 3).Value2 = "Los Lagos"
$ws.Cells.Item(58, 4).Value2 = 44498
$ws.Cells.Item(58, 5).Value2 = 10
$ws.Cells.Item(58, 6).Value2 = 100114014
$ws.Cells.Item(58, 7).Value2 = "Betarraga"
$ws.Cells.Item(58, 8).Value2 = "Sin especificar"
$ws.Cells.Item(58, 9).Value2 = "Primera"
$ws.Cells.Item(58, 10).Value2 = 1400
$ws.Cells.Item(58, 11).Value2 = 900
$ws.Cells.Item(58, 12).Value2 = 1000
$ws.Cells.Item(58, 13).Value2 = 950
$ws.Cells.Item(58, 14).Value2 = "`$/paquete 5 unidades"
$ws.Cells.Item(58, 15).Value2 = "Región del Maule"
$ws.Cells.Item(58, 16).Value2 = 190
$ws.Cells.Item(58, 17).Value2 = 5
$ws.Cells.Item(58, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date style ("s=2") used by every
# other date cell in column D.
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
